$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns (D:H) after the existing Rank/Team/ExpPoints columns so the
# new metric columns (WIN/TOP4/TOP5/TOP6/RELEGATION) sit between Team and the
# (relocated) ExpPoints column.
$ws.Range("D1:H1").EntireColumn.Insert()

# --- Header row -------------------------------------------------------
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"
$ws.Range("H1").Value = "ExpPoints"

# --- Re-order three clubs (rows 12-14) --------------------------------
$ws.Range("B12").Value = "Sevilla"
$ws.Range("B13").Value = "Valencia"
$ws.Range("B14").Value = "Espanyol"

# --- Recalculated ExpPoints values move from column C to column H; the
#     newly inserted WIN/TOP4/TOP5/TOP6/RELEGATION columns are left blank
#     for the upcoming Monte-Carlo simulation work. -------------------
$expPoints = @{
    2  = 88.63733836678766
    3  = 88.2838633875683
    4  = 70.07403845736171
    5  = 64.41179771828808
    6  = 61.25106960979502
    7  = 56.93389331042831
    8  = 53.1014170709494
    9  = 49.08106728033445
    10 = 48.82887536874397
    11 = 48.24543883537496
    12 = 47.52047920846063
    13 = 46.78673429747457
    14 = 46.17704639251744
    15 = 46.17383287980892
    16 = 40.19809224319757
    17 = 39.72843177731009
    18 = 38.38677328108188
    19 = 35.56210540595981
    20 = 33.05978326196885
    21 = 29.91713158274803
}

foreach ($row in 2..21) {
    $ws.Range("C$row").Value = ""
    $ws.Range("D$row").Value = ""
    $ws.Range("E$row").Value = ""
    $ws.Range("F$row").Value = ""
    $ws.Range("G$row").Value = ""
    $ws.Range("H$row").Value = $expPoints[$row]
}
